# New crime data collected - weekly CompStat update (cs-en-us-084pct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/issue number and reporting week dates ---
$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Row 16: Robbery ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -30.769230769230
$ws.Range("I16").Value = 135
$ws.Range("J16").Value = 138
$ws.Range("K16").Value = -2.173913043478
$ws.Range("L16").Value = 25
$ws.Range("M16").Value = -25.414364640884
$ws.Range("N16").Value = -89.647239263803

# --- Row 17: Fel. Assault ---
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -21.428571428571
$ws.Range("I17").Value = 172
$ws.Range("J17").Value = 170
$ws.Range("K17").Value = 1.176470588235
$ws.Range("L17").Value = 10.967741935483
$ws.Range("M17").Value = 53.571428571428
$ws.Range("N17").Value = -59.718969555035

# --- Row 18: Burglary ---
$ws.Range("C18").Value = "0"
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -66.666666666666
$ws.Range("J18").Value = 158
$ws.Range("K18").Value = 10.126582278481
$ws.Range("M18").Value = 47.457627118644
$ws.Range("N18").Value = -78.358208955223

# --- Row 19: Gr. Larceny ---
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 70
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -4.081632653061
$ws.Range("I19").Value = 628
$ws.Range("J19").Value = 529
$ws.Range("K19").Value = 18.714555765595
$ws.Range("L19").Value = 53.170731707317
$ws.Range("M19").Value = 22.65625
$ws.Range("N19").Value = -43.978590544157

# --- Row 20: G.L.A. ---
$ws.Range("C20").Value = "0"
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 66
$ws.Range("K20").Value = 46.666666666666
$ws.Range("L20").Value = 29.411764705882
$ws.Range("M20").Value = 29.411764705882
$ws.Range("N20").Value = -90.909090909090

# --- Row 21: TOTAL ---
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -8.333333333333
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 106
$ws.Range("H21").Value = -25.471698113207
$ws.Range("I21").Value = 1186
$ws.Range("J21").Value = 1051
$ws.Range("K21").Value = 12.844909609895
$ws.Range("L21").Value = 30.043859649122
$ws.Range("M21").Value = 21.144024514811
$ws.Range("N21").Value = -73.137032842582

# --- Row 22: Transit ---
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 57
$ws.Range("K22").Value = -1.754385964912
$ws.Range("L22").Value = -6.666666666666
$ws.Range("M22").Value = -22.222222222222

# --- Row 23: Housing ---
$ws.Range("C23").Value = "0"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = -20
$ws.Range("J23").Value = 52
$ws.Range("K23").Value = -25

# --- Row 24: Petit Larceny ---
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 17.391304347826
$ws.Range("F24").Value = 133
$ws.Range("G24").Value = 110
$ws.Range("H24").Value = 20.909090909090
$ws.Range("I24").Value = 1895
$ws.Range("J24").Value = 1328
$ws.Range("K24").Value = 42.695783132530
$ws.Range("L24").Value = 99.683877766069
$ws.Range("M24").Value = 27.352150537634

# --- Row 25: Misd. Assault ---
$ws.Range("C25").Value = "0"
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -19.230769230769
$ws.Range("J25").Value = 271
$ws.Range("K25").Value = 30.996309963099
$ws.Range("L25").Value = 37.596899224806
$ws.Range("M25").Value = 0.282485875706

# --- Row 26: UCR Rape* ---
$ws.Range("C26").Value = "0"
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0

# --- Row 27: Other Sex Crimes ---
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -70
$ws.Range("I27").Value = 64
$ws.Range("J27").Value = 69
$ws.Range("K27").Value = -7.246376811594
$ws.Range("L27").Value = 60
